# "complete session runs for removing records"
#
# The "Remove Incomplete Records" sheet modeled what the session-run table
# looks like once rows with missing runs are dropped. Rows 11-13 already had
# real numbers; rows 14-21 were still placeholder blanks. This fills in the
# remaining session-run values for rows 14-21 (columns C:L), which lets the
# MIN/MAX/MEAN/MEDIAN formulas in M:P resolve to real numbers instead of
# 0 / #DIV/0! / #NUM!.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Remove Incomplete Records")

$data = @(
    @(61.39,63.99,61.82,65.29,63.99,63.34,62.91,62.26,62.04,51.84),
    @(56.2,65,61,64.8,55.8,65.8,65.8,63.6,62,61),
    @(59.85,64.13,63.38,60.41,60.78,57.81,63.57,64.68,59.11,64.13),
    @(60.59,60.07,64.93,60.42,61.28,65.28,67.71,64.24,59.55,63.54),
    @(56.91,66.34,61.79,63.9,64.23,60.81,62.6,63.09,63.58,65.04),
    @(62.94,60.95,63.4,63.25,63.86,59.57,59.72,62.94,62.02,58.65),
    @(51.59,60.12,62.57,54.62,61.42,56.5,64.31,48.7,64.88,49.71),
    @(61.64,61.78,63.15,62.74,57.53,62.05,40.55,63.97,63.01,61.92)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 14 + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = 3 + $j   # column C = 3
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}

# The author left the cursor on M29 (below the filled table) on this sheet.
$ws.Activate()
$ws.Range("M29").Select()

# Active sheet moved from "Replace With Gradient" to "Replace With Mean".
$meanWs = $wb.Worksheets.Item("Replace With Mean")
$meanWs.Activate()
